# BookCover.xlsx update:
#  - insert a new "Date Published" column (new column D), shifting
#    Description/License List/Copyright/Authorship/JSON Manifest/Source/
#    Authorship Resource one column to the right
#  - fill in the Date Published values for each row
#  - tweak the Label text in row 4 (C4)
#  - fix up hyperlinks so they point at the shifted cells
#  - restore column widths / selection to match the edited layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D ("Description"); existing D..J shift to E..K.
$ws.Columns("D:D").Insert()

# 2. New column D header + values ("Date Published").
$ws.Range("D1").Value = "Date Published"
$ws.Range("D2").Value = "1928"
$ws.Range("D3").Value = "1916"
$ws.Range("D4").Value = "1915"
$ws.Range("D5").Value = "1921"
$ws.Range("D6").Value = "1913"
$ws.Range("D7").Value = "1900"
$ws.Range("D8").Value = "1899"

# 3. Updated Label text for row 4.
$ws.Range("C4").Value = "Alice in Wonderland - a dramatization of Lewis Carroll's ""Alice's adventures in Wonderland"" and ""Through the looking glass"", 1915"

# 4. Column widths: C and D both become the same (wider) width, and the
#    rest keep their previous widths (now shifted one column right).
$ws.Columns("C:D").ColumnWidth = 119.83203125

# 5. Hyperlinks used to live in H/I; after the column insert they belong in
#    I/J. The engine does not shift hyperlink ranges automatically on
#    Insert(), so rebuild them explicitly (delete-all, then re-add in the
#    exact original order so relationship ids come out the same way).
$ws.Range("A1").Hyperlinks.Delete()

$hls = $ws.Hyperlinks
[void]$hls.Add($ws.Range("I2"), "https://www.loc.gov/item/28026941/manifest.json")
[void]$hls.Add($ws.Range("J2"), "https://www.loc.gov/resource/gdcmassbookdig.alicesadventures00carr_25/?sp=1&st=single&r=-0.741,-0.078,2.482,1.56,0")
[void]$hls.Add($ws.Range("I3"), "https://www.loc.gov/item/16014724/manifest.json")
[void]$hls.Add($ws.Range("J3"), "https://www.loc.gov/item/16014724")
[void]$hls.Add($ws.Range("I4"), "https://www.loc.gov/item/16000724/manifest.json")
[void]$hls.Add($ws.Range("J4"), "https://www.loc.gov/resource/gdcmassbookdig.aliceinwonderlan00gers/?sp=11&r=-0.603,0.151,2.205,1.386,0")
[void]$hls.Add($ws.Range("I5"), "https://www.loc.gov/item/21027480/manifest.json")
[void]$hls.Add($ws.Range("J5"), "https://www.loc.gov/item/21027480/")
[void]$hls.Add($ws.Range("I6"), "https://www.loc.gov/item/42000114/manifest.json")
[void]$hls.Add($ws.Range("J6"), "https://www.loc.gov/item/42000114/")
[void]$hls.Add($ws.Range("B7"), "https://tile.loc.gov/image-services/iiif/public:gdcmassbookdig:throughlookinggl00carr_3:throughlookinggl00carr_3_0009/full/pct:100.0/0/default.jpg")
[void]$hls.Add($ws.Range("I7"), "https://www.loc.gov/item/00004842/manifest.json")
[void]$hls.Add($ws.Range("J7"), "https://www.loc.gov/item/00004842/")
[void]$hls.Add($ws.Range("B8"), "https://tile.loc.gov/image-services/iiif/public:gdcmassbookdig:throughlookinggl00carr_4:throughlookinggl00carr_4_0001/full/pct:100.0/0/default.jpg")
[void]$hls.Add($ws.Range("I8"), "https://www.loc.gov/item/00000848/manifest.json")
[void]$hls.Add($ws.Range("J8"), "https://www.loc.gov/item/00000848/")

# 6. Selection moves to D1 (matches the saved view after the edit).
$ws.Range("D1").Select()
